$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 3292
$ws.Range("J3").Value = 3417
$ws.Range("G4").Value = 1464
$ws.Range("J4").Value = 755
$ws.Range("J5").Value = 265
$ws.Range("J6").Value = 4042
$ws.Range("G7").Value = 24688
$ws.Range("J7").Value = 11771

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 46
$ws.Range("J7").Value = 136

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J4").Value = 9
$ws.Range("J6").Value = 40
$ws.Range("J7").Value = 129

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 372

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 89
$ws.Range("J3").Value = 173
$ws.Range("J6").Value = 120
$ws.Range("J7").Value = 421

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 88

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 88
$ws.Range("J3").Value = 93
$ws.Range("J7").Value = 309

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J5").Value = 32
$ws.Range("J7").Value = 357
$ws.Range("J8").Value = 764
$ws.Range("J11").Value = 178
$ws.Range("J15").Value = 136
$ws.Range("J16").Value = 33
$ws.Range("J18").Value = 116
$ws.Range("J19").Value = 367
$ws.Range("J20").Value = 248
$ws.Range("J21").Value = 21
$ws.Range("J30").Value = 52
$ws.Range("J31").Value = 88
$ws.Range("J33").Value = 517
$ws.Range("J36").Value = 171
$ws.Range("J37").Value = 372
$ws.Range("J42").Value = 467
$ws.Range("J43").Value = 107
$ws.Range("J45").Value = 16
$ws.Range("J46").Value = 43
$ws.Range("J48").Value = 116
$ws.Range("J49").Value = 75
$ws.Range("J50").Value = 67
$ws.Range("J51").Value = 157
$ws.Range("J53").Value = 114
$ws.Range("J54").Value = 225
$ws.Range("J55").Value = 148
$ws.Range("J57").Value = 53
$ws.Range("G63").Value = 221
$ws.Range("J63").Value = 61
$ws.Range("J65").Value = 309
$ws.Range("J67").Value = 421
$ws.Range("J76").Value = 166
$ws.Range("J77").Value = 99
$ws.Range("J78").Value = 155
$ws.Range("J79").Value = 350
$ws.Range("J83").Value = 273
$ws.Range("J85").Value = 538
$ws.Range("J88").Value = 121
$ws.Range("J89").Value = 136
$ws.Range("J90").Value = 138
$ws.Range("J91").Value = 135
$ws.Range("J92").Value = 37
$ws.Range("J94").Value = 103
$ws.Range("J96").Value = 129
$ws.Range("J97").Value = 73
$ws.Range("G101").Value = 24688
$ws.Range("J101").Value = 11771

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 80
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 273

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 141
$ws.Range("J4").Value = 26
$ws.Range("J5").Value = 19
$ws.Range("J7").Value = 517

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 59
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 225

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 208
$ws.Range("J3").Value = 231

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 89
$ws.Range("J7").Value = 367

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 55
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J6").Value = 89
$ws.Range("J7").Value = 166

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 132
$ws.Range("J3").Value = 202
$ws.Range("J4").Value = 41
$ws.Range("J6").Value = 152
$ws.Range("J7").Value = 538

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 97
$ws.Range("J6").Value = 231
$ws.Range("J7").Value = 467

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 39
$ws.Range("J7").Value = 155

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J6").Value = 75
$ws.Range("J7").Value = 148

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 39
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 129
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 350

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 78
$ws.Range("J7").Value = 248

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J2").Value = 35
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 171

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J3").Value = 20

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 43
$ws.Range("J7").Value = 136

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J2").Value = 19

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 66
$ws.Range("J3").Value = 36
$ws.Range("J5").Value = 4
$ws.Range("J6").Value = 60
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 121

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 227
$ws.Range("J3").Value = 239
$ws.Range("J6").Value = 235
$ws.Range("J7").Value = 764

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J2").Value = 45
$ws.Range("J7").Value = 138

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 107

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J2").Value = 26
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 114

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 34
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("J6").Value = 6
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J6").Value = 116
$ws.Range("J7").Value = 357

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J2").Value = 5
$ws.Range("J7").Value = 33
